$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.341.71'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7201'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.71'
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07999'
$ws.Range('E8').Value = '  +2.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3150'
$ws.Range('E9').Value = '  +1.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.98'
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08212'
$ws.Range('E11').Value = '  -2.21%  '
$ws.Range('D12').Value = '1.875.46'
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '94.73'
$ws.Range('E13').Value = '  +3.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.228'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7125'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.422'
$ws.Range('E16').Value = '  +5.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008498'
$ws.Range('E17').Value = '  +4.05%  '
$ws.Range('D18').Value = '29.337.53'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.59'
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.757'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.1595'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.038'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.52'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.405'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.309'
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.200'
$ws.Range('E31').Value = '  -7.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05364'
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('E33').Value = '  -0.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7616'
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.708'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01873'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '1.275.90'
$ws.Range('E38').Value = '  +3.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.751'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.449'
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '112.83'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9104'
$ws.Range('E42').Value = '  +2.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '74.09'
$ws.Range('E43').Value = '  +2.31%  '
$ws.Range('E44').Value = '  +7.27%  '
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').Value = '2.025.33'
$ws.Range('E46').Value = '  +0.35%  '
$ws.Range('E47').Value = '  +0.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.795'
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.507'
$ws.Range('E49').Value = '  +0.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4343'
$ws.Range('E50').Value = '  +0.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.089'
$ws.Range('E51').Value = '  +0.24%  '
